# Apply the coin-price / 1h-volume refresh from the GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.757.41"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.475.37"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.13"
$ws.Range("E5").Value = "  +1.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.19"
$ws.Range("E6").Value = "  +2.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.552"
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0880"
$ws.Range("E10").Value = "  +11.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "33.40"
$ws.Range("E11").Value = "  +2.96%  "
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.859.21"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.93"
$ws.Range("E14").Value = "  +1.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.64"
$ws.Range("E15").Value = "  -1.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.484.92"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.799"
$ws.Range("E17").Value = "  +3.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.730.21"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0952"
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.33"
$ws.Range("E22").Value = "  +2.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.15"
$ws.Range("E23").Value = "  +2.08%  "
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("E25").Value = "  +2.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.14"
$ws.Range("E27").Value = "  +2.52%  "
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.77"
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.07"
$ws.Range("E30").Value = "  +5.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.73"
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.51"
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0767"
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.57"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.44"
$ws.Range("E36").Value = "  +1.76%  "
$ws.Range("E37").Value = "  +5.68%  "
$ws.Range("E38").Value = "  +2.08%  "
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.104"
$ws.Range("E40").Value = "  +1.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.03"
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.52"
$ws.Range("E42").Value = "  +7.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.002.55"
$ws.Range("E43").Value = "  +3.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.13"
$ws.Range("E44").Value = "  +2.77%  "
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.98"
$ws.Range("E46").Value = "  +3.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.52"
$ws.Range("E47").Value = "  +5.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.719.23"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.74"
$ws.Range("E49").Value = "  +7.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.98"
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "67.46"
$ws.Range("E51").Value = "  +1.03%  "
